# Actualización automática 2025-06-02 14:03:43
# Adds a "PRESUPUESTO" (budget) column G to the "VENTA MENSUAL" sheet,
# mirroring the look (styles/column width) of the existing "junio" column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Clone formatting (header style, data style, total-row style) from column F.
$ws.Range("F1:F22").Copy()
$ws.Range("G1:G22").PasteSpecial(-4122)

# New column width of 17 characters (matches the "width=17" stored by Excel;
# ColumnWidth includes a built-in ~0.8333 padding offset).
$ws.Columns.Item(7).ColumnWidth = 16.166666666666668

# Header
$ws.Cells.Item(1, 7).Value = "PRESUPUESTO"

# Budget values for rows 2..21 (one per client), in row order; row 22 (last)
# is the column total.
$values = @(0, 1000, 0, 3000, 5000, 0, 6000, 0, 7000, 0, 6000, 6000, 1000, 400, 6500, 0, 4000, 0, 500, 4000, 50400)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $values[$i]
}
